# Update "想去人数" (interested-count) figures on the 展览 and 全部类型 sheets
# to match the latest scrape snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 297
$wsExhibit.Range("F6").Value = 152
$wsExhibit.Range("F7").Value = 303
$wsExhibit.Range("F8").Value = 219
$wsExhibit.Range("F9").Value = 2085
$wsExhibit.Range("F11").Value = 5036
$wsExhibit.Range("F12").Value = 104

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 297
$wsAll.Range("F8").Value = 152
$wsAll.Range("F9").Value = 303
$wsAll.Range("F10").Value = 219
$wsAll.Range("F13").Value = 2085
$wsAll.Range("F15").Value = 5036
$wsAll.Range("F16").Value = 104
